$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: dataset file name label
$ws.Range("A1").Value = "SemTab Round 2 - 211.csv & 212.csv"

# Row 3: president
$ws.Range("A3").Value = 'ResultatDataset1 "0" "1" URI: http://dbpedia.org/ontology/president'
$ws.Range("D3").Value = "http://dbpedia.org/ontology/president"

# Row 4: staff
$ws.Range("A4").Value = 'ResultatDataset1 "1" "1" URI: http://dbpedia.org/ontology/staff'
$ws.Range("D4").Value = "http://dbpedia.org/ontology/staff"

# Row 5: facultySize
$ws.Range("A5").Value = 'ResultatDataset1 "2" "1" URI: http://dbpedia.org/ontology/facultySize'
$ws.Range("D5").Value = "http://dbpedia.org/ontology/facultySize"

# Row 6: city
$ws.Range("A6").Value = 'ResultatDataset1 "3" "1" URI: http://dbpedia.org/ontology/city'
$ws.Range("D6").Value = "http://dbpedia.org/ontology/city"

# Row 7 used to be a data row (genus); now clear A7/D7 and turn it into SUM row over B3:B6 etc.
$ws.Range("A7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("B7").Formula = "=SUM(B3:B6)"
$ws.Range("C7").Formula = "=SUM(C3:C6)"
$ws.Range("E7").Formula = "=SUM(E3:E6)"

# Row 8 used to be the SUM row; now becomes Precision row
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("A8").Value = "Précision "
$ws.Range("B8").Formula = "=B7/C7"
$ws.Range("B8").NumberFormat = "0.00%"

# Row 9 used to be Precision row; now becomes Recall row
$ws.Range("A9").Value = "Recall "
$ws.Range("B9").Formula = "=C7/E7"

# Row 10 used to be Recall row; now becomes F1_Score row
$ws.Range("A10").Value = "F1_Score"
$ws.Range("B10").Formula = "=(2*B8*B9)/(B8+B9)"

# Row 11 used to be F1_Score row; now delete entirely
$ws.Rows(11).Delete()

# Update the active selection shown in the sheet view
$ws.Range("D11").Select()

$excel.CalculateFull()
